$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an exact text value into a cell without Excel re-interpreting
# numeric-/percent-looking strings as numbers (which would change both the
# stored type and silently lose precision, e.g. "0.002000" -> 0.002).
# A scratch cell well outside the used range is formatted as Text, loaded
# with the literal string, then copied in via PasteSpecial(xlPasteValues)
# so the destination cell keeps its original (default) style.
$scratch = $ws.Range("Z1")
function Set-TextValue([object]$range, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '246.51'
Set-TextValue $ws.Range('E2') '0.59%'
Set-TextValue $ws.Range('D3') '29.63'
Set-TextValue $ws.Range('E3') '9.78%'
Set-TextValue $ws.Range('D4') '5.163'
Set-TextValue $ws.Range('E4') '1.84%'
Set-TextValue $ws.Range('E5') '0.30%'
Set-TextValue $ws.Range('E6') '1.79%'
Set-TextValue $ws.Range('D7') '3.074'
Set-TextValue $ws.Range('E7') '2.38%'
Set-TextValue $ws.Range('D8') '0.8581'
Set-TextValue $ws.Range('E8') '4.69%'
Set-TextValue $ws.Range('D9') '0.8712'
Set-TextValue $ws.Range('E9') '3.70%'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D10') '0.01024'
Set-TextValue $ws.Range('E10') '1,613.53%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1364'
Set-TextValue $ws.Range('E11') '2.76%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.07084'
Set-TextValue $ws.Range('E12') '2.57%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.02921'
Set-TextValue $ws.Range('E13') '2.44%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.09392'
Set-TextValue $ws.Range('E14') '-0.05%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001516'
Set-TextValue $ws.Range('E15') '0.00%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D16') '0.04183'
Set-TextValue $ws.Range('E16') '2.39%'
Set-TextValue $ws.Range('D17') '0.006047'
Set-TextValue $ws.Range('E17') '-0.98%'
Set-TextValue $ws.Range('D19') '3.486'
Set-TextValue $ws.Range('E19') '-0.66%'
Set-TextValue $ws.Range('D20') '2.272'
Set-TextValue $ws.Range('E20') '-1.89%'
Set-TextValue $ws.Range('E21') '-0.09%'
Set-TextValue $ws.Range('D22') '0.03303'
Set-TextValue $ws.Range('E22') '3.37%'
Set-TextValue $ws.Range('D23') '0.1302'
Set-TextValue $ws.Range('E23') '0.38%'
Set-TextValue $ws.Range('D24') '3.476'
Set-TextValue $ws.Range('E24') '-2.43%'
Set-TextValue $ws.Range('E25') '0.52%'
Set-TextValue $ws.Range('D26') '0.005027'
Set-TextValue $ws.Range('E26') '26.81%'
Set-TextValue $ws.Range('D27') '0.001222'
Set-TextValue $ws.Range('E27') '0.28%'
Set-TextValue $ws.Range('E28') '23.58%'
Set-TextValue $ws.Range('D40') '0.03745'
Set-TextValue $ws.Range('D41') '0.005751'
Set-TextValue $ws.Range('E41') '4.64%'
Set-TextValue $ws.Range('D42') '0.1069'
Set-TextValue $ws.Range('E42') '1.39%'
Set-TextValue $ws.Range('D43') '0.002000'
Set-TextValue $ws.Range('E43') '11.20%'
Set-TextValue $ws.Range('D44') '0.009960'
Set-TextValue $ws.Range('E44') '6.03%'
Set-TextValue $ws.Range('D45') '0.00005212'
Set-TextValue $ws.Range('E45') '0.10%'
Set-TextValue $ws.Range('D47') '0.05801'
Set-TextValue $ws.Range('E47') '-42.81%'
Set-TextValue $ws.Range('D48') '0.002562'
Set-TextValue $ws.Range('E48') '-1.30%'

$scratch.Clear()
$excel.Application.CutCopyMode = $false

